$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row (row 3) with the latest 3R Games SA quote.
# Values are entered as plain text (matching the existing rows, which
# store numeric-looking figures and dates as text too), so numeric- and
# date-looking values are prefixed with a leading apostrophe to force
# text entry instead of Excel's automatic number/date detection.
$ws.Range("A3").Value = "3R Games SA"
$ws.Range("B3").Value = "0,88"
$ws.Range("C3").Value = "-1,13"
$ws.Range("D3").Value = "0,89"
$ws.Range("E3").Value = "0,88"
$ws.Range("F3").Value = "0,89"
$ws.Range("G3").Value = "'901"
$ws.Range("H3").Value = "'794"
$ws.Range("I3").Value = "'2021-07-02"
$ws.Range("J3").Value = "'2021-07-02"

# Drop the "quote prefix" formatting that entering text this way leaves
# behind, so the new cells don't pick up a style the original rows don't have.
$ws.Range("G3:J3").ClearFormats()
